$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D3").Value = 10.54
$ws.Range("E3").Value = 10.29

$ws.Range("C4").Value = 9.460000000000001
$ws.Range("E4").Value = 10.08
$ws.Range("F4").Value = 10.06

$ws.Range("C5").Value = 9.710000000000001
$ws.Range("D5").Value = 9.92
$ws.Range("F5").Value = 10.3
$ws.Range("G5").Value = 9.17
$ws.Range("J5").Value = 6.17

$ws.Range("D6").Value = 9.94
$ws.Range("E6").Value = 9.699999999999999
$ws.Range("G6").Value = 10.42

$ws.Range("E7").Value = 10.83
$ws.Range("F7").Value = 9.58
$ws.Range("H7").Value = 9.74

$ws.Range("G8").Value = 10.26

$ws.Range("E10").Value = 13.83
